$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q4" worksheet right after "总计" (i.e. before
#    the existing "2022-Q3" tab), and fill it with the quarterly fund-holder
#    breakdown.
# ---------------------------------------------------------------------------
$anchor = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item(1))
$q4.Name = "2022-Q4"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$rows = @(
    @("011220","南方匠心优选股票A","29.64","93.51","5.79","1.7162",7),
    @("020005","国泰金马稳健混合A","10.33","92.84","5.73","0.5919",6),
    @("005123","南方优享分红灵活配置混合A","6.79","92.15","8.40","0.5704",5),
    @("001044","嘉实新消费股票","10.68","81.38","5.24","0.5596",8),
    @("010446","国泰金福三个月定期开放混合","8.71","93.48","6.05","0.5270",4),
    @("000574","宝盈新价值灵活配置混合A","7.64","88.96","4.75","0.3629",5),
    @("003715","宝盈消费主题灵活配置混合","5.61","83.13","5.77","0.3237",4),
    @("011645","国泰核心价值两年持有期股票A","5.45","92.75","5.14","0.2801",7),
    @("004357","南方智慧精选灵活配置混合","4.68","91.92","5.55","0.2597",9),
    @("006921","南方智诚混合","4.43","93.90","5.84","0.2587",7),
    @("008174","国泰蓝筹精选混合A","5.37","79.99","4.19","0.2250",7),
    @("009223","宝盈现代服务业混合A","3.46","88.59","5.01","0.1733",7),
    @("006587","南方优享分红灵活配置混合C","1.84","92.15","8.40","0.1546",5),
    @("005810","南方瑞祥一年定期开放灵活配置混合A","2.80","93.56","5.38","0.1506",9),
    @("007574","宝盈新价值灵活配置混合C","3.00","88.96","4.75","0.1425",5),
    @("011221","南方匠心优选股票C","1.76","93.51","5.79","0.1019",7),
    @("481017","工银量化策略混合A","2.42","92.03","3.77","0.0912",4),
    @("012308","国泰价值远见两年封闭运作混合A","1.62","93.18","5.48","0.0888",6),
    @("008175","国泰蓝筹精选混合C","1.71","79.99","4.19","0.0716",7),
    @("011384","南方远见回报股票A","1.04","93.53","5.64","0.0587",8),
    @("008303","宝盈龙头优选股票A","0.65","88.21","4.90","0.0318",8),
    @("009224","宝盈现代服务业混合C","0.41","88.59","5.01","0.0205",7),
    @("005811","南方瑞祥一年定期开放灵活配置混合C","0.32","93.56","5.38","0.0172",9),
    @("011727","工银瑞信聚瑞混合A","0.83","39.61","1.68","0.0139",8),
    @("011385","南方远见回报股票C","0.21","93.53","5.64","0.0118",8),
    @("011646","国泰核心价值两年持有期股票C","0.19","92.75","5.14","0.0098",7),
    @("006675","宝盈品牌消费股票A","0.18","87.69","5.01","0.0090",7),
    @("008304","宝盈龙头优选股票C","0.18","88.21","4.90","0.0088",8),
    @("001648","工银新价值灵活配置混合A","0.51","87.39","1.60","0.0082",8),
    @("006676","宝盈品牌消费股票C","0.15","87.69","5.01","0.0075",7),
    @("012309","国泰价值远见两年封闭运作混合C","0.10","93.18","5.48","0.0055",6),
    @("012237","工银新价值灵活配置混合C","0.03","87.39","1.60","0.0005",8),
    @("012241","工银量化策略混合C","0.01","92.03","3.77","0.0004",4),
    @("015589","国泰金马稳健混合C","0.00","92.84","5.73","0",6),
    @("011728","工银瑞信聚瑞混合C","0.00","39.61","1.68","0",8)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $q4.Cells.Item($r, 1).Value = $i
    $q4.Cells.Item($r, 2).Value = $row[0]
    $q4.Cells.Item($r, 3).Value = $row[1]
    $q4.Cells.Item($r, 4).Value = $row[2]
    $q4.Cells.Item($r, 5).Value = $row[3]
    $q4.Cells.Item($r, 6).Value = $row[4]
    $q4.Cells.Item($r, 7).Value = $row[5]
    $q4.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: push every existing row down by one,
#    insert the new 2022-Q4 totals at the top, and append the 2020-Q4 row
#    that now appears at the bottom of the (now one-row-longer) history.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$totalRows = @(
    @("2022-Q4", 35, 6.85),
    @("2022-Q3", 54, 10.54),
    @("2022-Q2", 27, 4.35),
    @("2022-Q1", 20, 2.07),
    @("2021-Q4", 24, 3.53),
    @("2021-Q3", 8, 0.47),
    @("2021-Q1", 1, 0.04),
    @("2020-Q4", 1, 0)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}

# Column A holds a plain running index (0..7); rows 2-8 already carry it
# with the correct style, only row 9 is brand new and needs both the value
# and the matching (bordered) formatting copied over from the row above it.
$total.Cells.Item(8, 1).Copy()
$total.Cells.Item(9, 1).PasteSpecial(-4122)
$total.Cells.Item(9, 1).Value = 7
